$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new shared string (description for the new row) ---
$newDescription = "Aloin suunitelemaana lisää ominaisuuksia. Implementoin uusia usercontrol ja niille omat navigaatio osiot. Lajitelin tiedostoja paremmin. Tein random number generaatorin. Tein Yksinkertaisen checksum työkalun"

# --- Move the current totals row (row 24) down to row 25 ---
# Copy formats of the totals row first so the new row keeps the same style ids
# (avoids creating duplicate cellXfs entries the way Rows.Insert() would).
$ws.Range("B24:D24").Copy()
$ws.Range("B25:D25").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B25").Value2 = "Yht"
$ws.Range("C25").Formula = "=SUM(C6:C24)"
$ws.Range("D25").Value2 = ""
$ws.Rows(25).RowHeight = 18.75

# --- Build the new data row 24 using the same look as the rows above it ---
$ws.Range("B23:D23").Copy()
$ws.Range("B24:D24").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B24").Value = 45345
$ws.Range("C24").Value2 = 9
$ws.Range("D24").Value2 = $newDescription

$ws.Rows(24).RowHeight = 93.75

# --- Update the selection to match the authored state ---
$ws.Range("G24").Select()

$wb.Save()
